$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.79769999999999
$ws.Range("E4").Value = 12.87540000000001
$ws.Range("A9").Value = -19.91999999999999
$ws.Range("E10").Value = 12.2022
$ws.Range("A18").Value = -23.02970000000002
$ws.Range("A20").Value = -22.15610000000002
$ws.Range("D21").Value = -7.426100000000003
